{"js": "// Replace each three-digit-by-one-digit multiplication problem with its\n// new value, in place, matching the exact old text so there is no\n// ambiguity. All 25 values are unique within the document.\nconst replacements = [\n  [\"135\u00d72=\", \"491\u00d76=\"],\n  [\"579\u00d72=\", \"134\u00d74=\"],\n  [\"902\u00d73=\", \"740\u00d78=\"],\n  [\"726\u00d72=\", \"686\u00d78=\"],\n  [\"455\u00d77=\", \"952\u00d79=\"],\n  [\"419\u00d79=\", \"550\u00d72=\"],\n  [\"193\u00d72=\", \"305\u00d73=\"],\n  [\"132\u00d72=\", \"273\u00d74=\"],\n  [\"232\u00d76=\", \"617\u00d78=\"],\n  [\"717\u00d75=\", \"750\u00d76=\"],\n  [\"826\u00d72=\", \"646\u00d78=\"],\n  [\"682\u00d72=\", \"580\u00d73=\"],\n  [\"185\u00d78=\", \"808\u00d75=\"],\n  [\"623\u00d74=\", \"840\u00d78=\"],\n  [\"976\u00d74=\", \"105\u00d74=\"],\n  [\"714\u00d76=\", \"597\u00d72=\"],\n  [\"486\u00d76=\", \"831\u00d76=\"],\n  [\"446\u00d73=\", \"635\u00d77=\"],\n  [\"838\u00d75=\", \"686\u00d79=\"],\n  [\"668\u00d79=\", \"707\u00d75=\"],\n  [\"958\u00d77=\", \"522\u00d77=\"],\n  [\"791\u00d75=\", \"838\u00d73=\"],\n  [\"898\u00d76=\", \"561\u00d77=\"],\n  [\"830\u00d77=\", \"679\u00d74=\"],\n  [\"540\u00d77=\", \"678\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication problem with its\n# new value, in place, matching the exact old text so there is no\n# ambiguity. All 25 values are unique within the document.\n$pairs = @(\n    @{Old=\"135\u00d72=\"; New=\"491\u00d76=\"},\n    @{Old=\"579\u00d72=\"; New=\"134\u00d74=\"},\n    @{Old=\"902\u00d73=\"; New=\"740\u00d78=\"},\n    @{Old=\"726\u00d72=\"; New=\"686\u00d78=\"},\n    @{Old=\"455\u00d77=\"; New=\"952\u00d79=\"},\n    @{Old=\"419\u00d79=\"; New=\"550\u00d72=\"},\n    @{Old=\"193\u00d72=\"; New=\"305\u00d73=\"},\n    @{Old=\"132\u00d72=\"; New=\"273\u00d74=\"},\n    @{Old=\"232\u00d76=\"; New=\"617\u00d78=\"},\n    @{Old=\"717\u00d75=\"; New=\"750\u00d76=\"},\n    @{Old=\"826\u00d72=\"; New=\"646\u00d78=\"},\n    @{Old=\"682\u00d72=\"; New=\"580\u00d73=\"},\n    @{Old=\"185\u00d78=\"; New=\"808\u00d75=\"},\n    @{Old=\"623\u00d74=\"; New=\"840\u00d78=\"},\n    @{Old=\"976\u00d74=\"; New=\"105\u00d74=\"},\n    @{Old=\"714\u00d76=\"; New=\"597\u00d72=\"},\n    @{Old=\"486\u00d76=\"; New=\"831\u00d76=\"},\n    @{Old=\"446\u00d73=\"; New=\"635\u00d77=\"},\n    @{Old=\"838\u00d75=\"; New=\"686\u00d79=\"},\n    @{Old=\"668\u00d79=\"; New=\"707\u00d75=\"},\n    @{Old=\"958\u00d77=\"; New=\"522\u00d77=\"},\n    @{Old=\"791\u00d75=\"; New=\"838\u00d73=\"},\n    @{Old=\"898\u00d76=\"; New=\"561\u00d77=\"},\n    @{Old=\"830\u00d77=\"; New=\"679\u00d74=\"},\n    @{Old=\"540\u00d77=\"; New=\"678\u00d76=\"}\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
